$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 441
$ws1.Range("F7").Value = 1291
$ws1.Range("F8").Value = 489
$ws1.Range("F10").Value = 256
$ws1.Range("F17").Value = 58
$ws1.Range("F18").Value = 232
$ws1.Range("F19").Value = 1615
$ws1.Range("F21").Value = 260
$ws1.Range("F22").Value = 150
$ws1.Range("F23").Value = 1239
$ws1.Range("F27").Value = 1191
$ws1.Range("F30").Value = 2782
$ws1.Range("F31").Value = 1557
$ws1.Range("F33").Value = 92
$ws1.Range("F34").Value = 586
$ws1.Range("F36").Value = 1637
$ws1.Range("F37").Value = 866
$ws1.Range("F38").Value = 1672
$ws1.Range("F40").Value = 89
$ws1.Range("F41").Value = 821
$ws1.Range("F42").Value = 19
$ws1.Range("F43").Value = 777
$ws1.Range("F44").Value = 759
$ws1.Range("F45").Value = 959
$ws1.Range("F46").Value = 416
$ws1.Range("F47").Value = 3296

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 184
$ws2.Range("F15").Value = 771

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 441
$ws4.Range("F9").Value = 1291
$ws4.Range("F10").Value = 489
$ws4.Range("F11").Value = 256
$ws4.Range("F18").Value = 1615
$ws4.Range("F20").Value = 150
$ws4.Range("F21").Value = 1240
$ws4.Range("F27").Value = 1191
$ws4.Range("F28").Value = 2782
$ws4.Range("F29").Value = 1557
$ws4.Range("F31").Value = 92
$ws4.Range("F32").Value = 771
$ws4.Range("F36").Value = 586
$ws4.Range("F37").Value = 1637
$ws4.Range("F40").Value = 866
$ws4.Range("F41").Value = 1672
$ws4.Range("F42").Value = 821
$ws4.Range("F43").Value = 777
$ws4.Range("F44").Value = 759
$ws4.Range("F45").Value = 959
$ws4.Range("F46").Value = 416
$ws4.Range("F48").Value = 3296
